$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "V1"
$ws.Range("I1").Value = "V2"
$ws.Range("J1").Value = "V3"
$ws.Range("K1").Value = "V4"
$ws.Range("L1").Value = "V5"

$ws.Range("L1").Select()
